# Append the 2025-09-11 daily totals for both charging stations
# (four-square-lawn station / 四方坪站 and gaoling station / 高岭站)
# as two new rows at the bottom of the data table, then move the
# active-cell selection down to reflect the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84 - 四方坪站 (site code matches the string already used in row 80 / B4 etc.)
$ws.Cells.Item(84, 1).Value = 45911
$ws.Cells.Item(84, 2).Value = "四方坪站"
$ws.Cells.Item(84, 3).Value = 11156.15
$ws.Cells.Item(84, 4).Value = 8997.87
$ws.Cells.Item(84, 5).Value = 3931.89
$ws.Cells.Item(84, 6).Value = 458

# Row 85 - 高岭站
$ws.Cells.Item(85, 1).Value = 45911
$ws.Cells.Item(85, 2).Value = "高岭站"
$ws.Cells.Item(85, 3).Value = 4379.29
$ws.Cells.Item(85, 4).Value = 3461.97
$ws.Cells.Item(85, 5).Value = 1104.8
$ws.Cells.Item(85, 6).Value = 160

# Update the selected cell, mirroring the workbook's recorded selection move
$ws.Range("H79").Select() | Out-Null
